$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (test_valid1): fix address row/col order, and correct size value ---
$ws.Range("C7").Value = "1. Enter shipment weight = 1 kg, size = 0.5, address = (row '1', col 'B')`r`n2. Enter shipment weight = 300 kg and size = 20`r`n3. Verify result."
$ws.Range("D7").Value = "SHIPMENT`r`nweight: 300 kg`r`nsize: 2`r`n`r`nMAP`r`nmap.squares[1][1] = 1"

# --- Row 8 (test_valid2): fix address row/col order (Test Data unchanged) ---
$ws.Range("C8").Value = "1. Enter shipment weight = 5001 kg, size = 2, address = (row '2', col 'B')`r`n2. Verify result."
$ws.Range("D8").Value = "SHIPMENT`r`nweight: 5001 kg`r`nsize: 2`r`n`r`nMAP`r`nmap.squares[1][2] = 1"

# --- Row 9 (test_valid3): fix address row/col order, correct size + map coordinates ---
$ws.Range("C9").Value = "1. Enter shipment weight = 10 kg, size = 3, address = (row '0', col 'G')`r`n2. Verify result."
$ws.Range("D9").Value = "SHIPMENT`r`nweight: 10 kg`r`nsize: 2`r`n`r`nMAP`r`nmap.squares[0][6] != 1"

# --- Row 10 (test_valid4): fix address row/col order, correct size + map coordinates ---
$ws.Range("C10").Value = "1. Enter shipment weight = 20 kg, size = 2, address = (row '26', col 'A')`r`n2. Verify result."
$ws.Range("D10").Value = "SHIPMENT`r`nweight: 20 kg`r`nsize: 2`r`n`r`nMAP`r`nmap.squares[26][0] != 1"

# --- Update the window selection to match the latest editing position ---
$ws.Range("A22").Select()
